$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cBold = $ws.Range("A500")
$cBold.Value = "bold-test"
$cBold.Style = "Normal"
$cBold.Font.Bold = $true

$st = $wb.Styles.Add("MenloStyle")
$st.Font.Name = "Menlo"
$st.Font.Size = 12
$cMenlo = $ws.Range("A501")
$cMenlo.Value = "menlo-test"
$cMenlo.Style = "MenloStyle"
$wb.Styles("MenloStyle").Delete()
